$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "PETUNJUK" instruction row (row 1); the rows below shift up.
$ws.Rows("1").Delete()

# Row 2 (was row 3) used to hold the placeholder "[1]" / "Nama Kelas 1".
# Replace with a real numeric index and the proper uppercase class name.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "NAMA KELAS 1"

# Re-align that row to the right (replacing the old centered placeholder look).
$ws.Range("A2:B2").HorizontalAlignment = -4152
$ws.Range("A2:B2").VerticalAlignment = -4107

# Restore the selection left behind by the editing session.
$ws.Range("B3").Select() | Out-Null
